# Project DesignFirst save: update Rules sheet cell C10 (the R20 rule's
# "Integer min" threshold) from 18 to 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Rules")

$ws.Range("C10").Value = 100

$wb.Save()
